$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 6 (the time-series data shrank from 5 data rows to 4)
$ws.Rows.Item(6).Delete()

# 2) Overwrite rows 2-5 with the refreshed measurement values
# Row 2
$ws.Cells.Item(2, 1).Value = 45071.50694444445
$ws.Cells.Item(2, 2).Value = 12.77
$ws.Cells.Item(2, 3).Value = 8.425000000000001
$ws.Cells.Item(2, 4).Value = 3.488
$ws.Cells.Item(2, 5).Value = 27.776
$ws.Cells.Item(2, 6).Value = 20.765
$ws.Cells.Item(2, 7).Value = 9.82
$ws.Cells.Item(2, 8).Value = 29.191
$ws.Cells.Item(2, 9).Value = 15.707
$ws.Cells.Item(2, 10).Value = 6.218
$ws.Cells.Item(2, 11).Value = 9.112
$ws.Cells.Item(2, 12).Value = 10.928
$ws.Cells.Item(2, 13).Value = 11.76
$ws.Cells.Item(2, 14).Value = 3.255
$ws.Cells.Item(2, 15).Value = 10.151
$ws.Cells.Item(2, 16).Value = 13.871
$ws.Cells.Item(2, 17).Value = 9.122999999999999
$ws.Cells.Item(2, 18).Value = 2.77
$ws.Cells.Item(2, 19).Value = 1.592
$ws.Cells.Item(2, 20).Value = 147.125
$ws.Cells.Item(2, 21).Value = 28.067
$ws.Cells.Item(2, 22).Value = 9.369999999999999
$ws.Cells.Item(2, 23).Value = 18.004
$ws.Cells.Item(2, 24).Value = 9.226000000000001
$ws.Cells.Item(2, 25).Value = 2.728
$ws.Cells.Item(2, 26).Value = 15.858
$ws.Cells.Item(2, 27).Value = 8.276
$ws.Cells.Item(2, 28).Value = 7.652
$ws.Cells.Item(2, 29).Value = 9.069000000000001
$ws.Cells.Item(2, 30).Value = 11.378
$ws.Cells.Item(2, 31).Value = 2.869
$ws.Cells.Item(2, 32).Value = 26.432
$ws.Cells.Item(2, 33).Value = 4.957
$ws.Cells.Item(2, 34).Value = 11.713

# Row 3
$ws.Cells.Item(3, 1).Value = 45071.51388888889
$ws.Cells.Item(3, 2).Value = 11.856
$ws.Cells.Item(3, 3).Value = 8.467000000000001
$ws.Cells.Item(3, 4).Value = 1.608
$ws.Cells.Item(3, 5).Value = 26.068
$ws.Cells.Item(3, 6).Value = 20.377
$ws.Cells.Item(3, 7).Value = 9.208
$ws.Cells.Item(3, 8).Value = 35.568
$ws.Cells.Item(3, 9).Value = 14.543
$ws.Cells.Item(3, 10).Value = 6.178
$ws.Cells.Item(3, 11).Value = 8.929
$ws.Cells.Item(3, 12).Value = 10.395
$ws.Cells.Item(3, 13).Value = 11.155
$ws.Cells.Item(3, 14).Value = 3.02
$ws.Cells.Item(3, 15).Value = 9.398999999999999
$ws.Cells.Item(3, 16).Value = 13.144
$ws.Cells.Item(3, 17).Value = 8.301
$ws.Cells.Item(3, 18).Value = 1.343
$ws.Cells.Item(3, 19).Value = 0.88
$ws.Cells.Item(3, 20).Value = 135.739
$ws.Cells.Item(3, 21).Value = 26.303
$ws.Cells.Item(3, 22).Value = 8.676
$ws.Cells.Item(3, 23).Value = 17.252
$ws.Cells.Item(3, 24).Value = 8.962999999999999
$ws.Cells.Item(3, 25).Value = 1.822
$ws.Cells.Item(3, 26).Value = 17.744
$ws.Cells.Item(3, 27).Value = 7.663
$ws.Cells.Item(3, 28).Value = 6.986
$ws.Cells.Item(3, 29).Value = 8.227
$ws.Cells.Item(3, 30).Value = 10.843
$ws.Cells.Item(3, 31).Value = 1.165
$ws.Cells.Item(3, 32).Value = 32.754
$ws.Cells.Item(3, 33).Value = 4.704
$ws.Cells.Item(3, 34).Value = 10.847

# Row 4
$ws.Cells.Item(4, 1).Value = 45071.52083333334
$ws.Cells.Item(4, 2).Value = 20.533
$ws.Cells.Item(4, 3).Value = 15.183
$ws.Cells.Item(4, 4).Value = 1.442
$ws.Cells.Item(4, 5).Value = 44.922
$ws.Cells.Item(4, 6).Value = 36.329
$ws.Cells.Item(4, 7).Value = 16.079
$ws.Cells.Item(4, 8).Value = 61.097
$ws.Cells.Item(4, 9).Value = 25.014
$ws.Cells.Item(4, 10).Value = 11.01
$ws.Cells.Item(4, 11).Value = 16.2
$ws.Cells.Item(4, 12).Value = 17.996
$ws.Cells.Item(4, 13).Value = 19.169
$ws.Cells.Item(4, 14).Value = 5.193
$ws.Cells.Item(4, 15).Value = 16.166
$ws.Cells.Item(4, 16).Value = 22.917
$ws.Cells.Item(4, 17).Value = 13.767
$ws.Cells.Item(4, 18).Value = 0.998
$ws.Cells.Item(4, 19).Value = 0.9350000000000001
$ws.Cells.Item(4, 20).Value = 238.75
$ws.Cells.Item(4, 21).Value = 45.175
$ws.Cells.Item(4, 22).Value = 14.922
$ws.Cells.Item(4, 23).Value = 30.214
$ws.Cells.Item(4, 24).Value = 15.825
$ws.Cells.Item(4, 25).Value = 2.51
$ws.Cells.Item(4, 26).Value = 30.129
$ws.Cells.Item(4, 27).Value = 13.181
$ws.Cells.Item(4, 28).Value = 11.751
$ws.Cells.Item(4, 29).Value = 13.824
$ws.Cells.Item(4, 30).Value = 18.87
$ws.Cells.Item(4, 31).Value = 0.733
$ws.Cells.Item(4, 32).Value = 55.657
$ws.Cells.Item(4, 33).Value = 8.332000000000001
$ws.Cells.Item(4, 34).Value = 18.656

# Row 5
$ws.Cells.Item(5, 1).Value = 45071.52777777778
$ws.Cells.Item(5, 2).Value = 9.5
$ws.Cells.Item(5, 3).Value = 6.98
$ws.Cells.Item(5, 4).Value = 0.87
$ws.Cells.Item(5, 5).Value = 20.9
$ws.Cells.Item(5, 6).Value = 16.65
$ws.Cells.Item(5, 7).Value = 7.42
$ws.Cells.Item(5, 8).Value = 32.82
$ws.Cells.Item(5, 9).Value = 11.63
$ws.Cells.Item(5, 10).Value = 5.08
$ws.Cells.Item(5, 11).Value = 7.36
$ws.Cells.Item(5, 12).Value = 8.369999999999999
$ws.Cells.Item(5, 13).Value = 8.960000000000001
$ws.Cells.Item(5, 14).Value = 2.42
$ws.Cells.Item(5, 15).Value = 7.52
$ws.Cells.Item(5, 16).Value = 10.63
$ws.Cells.Item(5, 17).Value = 6.52
$ws.Cells.Item(5, 18).Value = 0.7
$ws.Cells.Item(5, 19).Value = 0.52
$ws.Cells.Item(5, 20).Value = 107.14
$ws.Cells.Item(5, 21).Value = 21.17
$ws.Cells.Item(5, 22).Value = 6.94
$ws.Cells.Item(5, 23).Value = 14.05
$ws.Cells.Item(5, 24).Value = 7.31
$ws.Cells.Item(5, 25).Value = 1.28
$ws.Cells.Item(5, 26).Value = 15.56
$ws.Cells.Item(5, 27).Value = 6.13
$ws.Cells.Item(5, 28).Value = 5.53
$ws.Cells.Item(5, 29).Value = 6.49
$ws.Cells.Item(5, 30).Value = 8.76
$ws.Cells.Item(5, 31).Value = 0.54
$ws.Cells.Item(5, 32).Value = 30.16
$ws.Cells.Item(5, 33).Value = 3.82
$ws.Cells.Item(5, 34).Value = 8.68

# 3) Widen a subset of data columns by 1 character
# NOTE: Excel's ColumnWidth (COM, character units) is re-derived from an
# internal pixel width, which adds a constant ~0.85-character pad before it
# is re-quantized back into "characters" for the OOXML <col width> attribute.
# Subtracting 0.85 from the desired OOXML width reliably round-trips to the
# exact integer width we want (verified empirically against this runtime).
$ws.Columns.Item(2).ColumnWidth = 7.15
$ws.Columns.Item(3).ColumnWidth = 7.15
$ws.Columns.Item(5).ColumnWidth = 7.15
$ws.Columns.Item(6).ColumnWidth = 7.15
$ws.Columns.Item(7).ColumnWidth = 7.15
$ws.Columns.Item(8).ColumnWidth = 7.15
$ws.Columns.Item(9).ColumnWidth = 7.15
$ws.Columns.Item(12).ColumnWidth = 7.15
$ws.Columns.Item(13).ColumnWidth = 7.15
$ws.Columns.Item(15).ColumnWidth = 7.15
$ws.Columns.Item(16).ColumnWidth = 7.15
$ws.Columns.Item(17).ColumnWidth = 7.15
$ws.Columns.Item(20).ColumnWidth = 8.15
$ws.Columns.Item(21).ColumnWidth = 7.15
$ws.Columns.Item(22).ColumnWidth = 7.15
$ws.Columns.Item(23).ColumnWidth = 7.15
$ws.Columns.Item(24).ColumnWidth = 7.15
$ws.Columns.Item(26).ColumnWidth = 7.15
$ws.Columns.Item(27).ColumnWidth = 7.15
$ws.Columns.Item(28).ColumnWidth = 7.15
$ws.Columns.Item(29).ColumnWidth = 7.15
$ws.Columns.Item(30).ColumnWidth = 7.15
$ws.Columns.Item(32).ColumnWidth = 7.15
$ws.Columns.Item(34).ColumnWidth = 7.15
